# This workbook is a weekly time series of "Mora" (blackberry) wholesale
# prices at "Mercado Mayorista Lo Valledor de Santiago". The commit adds one
# new weekly observation. Since the sheet is sorted with the newest-looking
# entries interleaved, the new record lands at row 14, and every existing
# data row from (old) row 14 through (old) row 89 shifts down by one.
#
# Inserting a whole row (rather than rewriting every cell of every row)
# reproduces exactly that shift, including the date-number style (s="2" on
# column D) that Excel automatically propagates to the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 14; everything currently at row 14.. shifts to row 15..
$ws.Rows(14).Insert()

# Populate the newly inserted row 14 with the new observation.
$ws.Cells.Item(14, 1).Value = 6
$ws.Cells.Item(14, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(14, 3).Value = "Metropolitana"
$ws.Cells.Item(14, 4).Value = 44901
$ws.Cells.Item(14, 5).Value = 13
$ws.Cells.Item(14, 6).Value = "Fruta"
$ws.Cells.Item(14, 7).Value = 100101
$ws.Cells.Item(14, 8).Value = "Berries"
$ws.Cells.Item(14, 9).Value = 100101008
$ws.Cells.Item(14, 10).Value = "Mora"
$ws.Cells.Item(14, 11).Value = "Sin especificar"
$ws.Cells.Item(14, 12).Value = "Especial"
$ws.Cells.Item(14, 13).Value = 250
$ws.Cells.Item(14, 14).Value = 6000
$ws.Cells.Item(14, 15).Value = 6000
$ws.Cells.Item(14, 16).Value = 6000
$ws.Cells.Item(14, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(14, 18).Value = "Región del Maule"
$ws.Cells.Item(14, 19).Value = 3000
$ws.Cells.Item(14, 20).Value = 2
